$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 430.66666
$ws.Range("I2").Value = 95.333336
$ws.Range("J2").Value = 766
$ws.Range("K2").Value = 95.333336
$ws.Range("L2").Value = 766
$ws.Range("M2").Value = 17.666664
$ws.Range("N2").Value = -992
$ws.Range("H9").Value = 120.6
$ws.Range("I9").Value = 124.63158
$ws.Range("K9").Value = 124.63158
$ws.Range("M9").Value = 44.36842
$ws.Range("H17").Value = 1988.1936
$ws.Range("J17").Value = 2177.36
$ws.Range("L17").Value = 6532.08
$ws.Range("N17").Value = -6868.08
$ws.Range("H32").Value = 0
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 0
$ws.Range("H58").Value = 301.22223
$ws.Range("I58").Value = 245.125
$ws.Range("J58").Value = 750
$ws.Range("K58").Value = 735.375
$ws.Range("L58").Value = 2250
$ws.Range("M58").Value = -585.375
$ws.Range("N58").Value = -2550
$ws.Range("H123").Value = 8379999.5
$ws.Range("J123").Value = 69999.5
$ws.Range("L123").Value = 69999.5
$ws.Range("N123").Value = -79799.5
$ws.Range("H132").Value = 74595.09
$ws.Range("I132").Value = 83339
$ws.Range("K132").Value = 250017
$ws.Range("M132").Value = -247487
$ws.Range("H137").Value = 898969.5600000001
$ws.Range("I137").Value = 435041.3
$ws.Range("K137").Value = 1305123.9
$ws.Range("M137").Value = -1302573.9
$ws.Range("H138").Value = 1004749.06
$ws.Range("J138").Value = 1474421
$ws.Range("L138").Value = 4423263
$ws.Range("N138").Value = -4433543
$ws.Range("M32").ClearContents()
$ws.Range("N32").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1855.079
$ws.Range("I2").Value = 1767.3793
$ws.Range("J2").Value = 2137.6667
$ws.Range("K2").Value = 1767.3793
$ws.Range("L2").Value = 2137.6667
$ws.Range("M2").Value = -1654.3793
$ws.Range("N2").Value = -2363.6667
$ws.Range("H26").Value = 1978.5
$ws.Range("I26").Value = 1978.5
$ws.Range("K26").Value = 1978.5
$ws.Range("M26").Value = -1648.5
$ws.Range("H32").Value = 5687234
$ws.Range("I32").Value = 5887783.5
$ws.Range("K32").Value = 5887783.5
$ws.Range("M32").Value = -5887496.5
$ws.Range("H45").Value = 5483.0835
$ws.Range("I45").Value = 5483.0835
$ws.Range("K45").Value = 5483.0835
$ws.Range("M45").Value = -5106.0835
$ws.Range("H61").Value = 802631.75
$ws.Range("I61").Value = 1118184.6
$ws.Range("K61").Value = 1118184.6
$ws.Range("M61").Value = -1117972.6
$ws.Range("H74").Value = 2086024.1
$ws.Range("I74").Value = 2453394
$ws.Range("K74").Value = 2453394
$ws.Range("M74").Value = -2452520
$ws.Range("H76").Value = 9666
$ws.Range("J76").Value = 9666
$ws.Range("L76").Value = 9666
$ws.Range("N76").Value = -10342
$ws.Range("H77").Value = 2086024.1
$ws.Range("I77").Value = 2453394
$ws.Range("K77").Value = 12266970
$ws.Range("M77").Value = -12262602
$ws.Range("H79").Value = 9666
$ws.Range("J79").Value = 9666
$ws.Range("L79").Value = 9666
$ws.Range("N79").Value = -12006
$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("H116").Value = 1855.079
$ws.Range("I116").Value = 1767.3793
$ws.Range("J116").Value = 2137.6667
$ws.Range("K116").Value = 1767.3793
$ws.Range("L116").Value = 2137.6667
$ws.Range("M116").Value = 526.6206999999999
$ws.Range("N116").Value = -6725.6667
$ws.Range("H132").Value = 251436.97
$ws.Range("I132").Value = 417885.88
$ws.Range("K132").Value = 1253657.64
$ws.Range("M132").Value = -1251127.64
$ws.Range("H136").Value = 802631.75
$ws.Range("I136").Value = 1118184.6
$ws.Range("K136").Value = 3354553.8
$ws.Range("M136").Value = -3352003.8
$ws.Range("N114").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1855.079
$ws.Range("I3").Value = 1767.3793
$ws.Range("J3").Value = 2137.6667
$ws.Range("K3").Value = 1767.3793
$ws.Range("L3").Value = 2137.6667
$ws.Range("M3").Value = -1653.3793
$ws.Range("N3").Value = -2365.6667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7002.3267
$ws.Range("I31").Value = 1864.5
$ws.Range("J31").Value = 9493.394
$ws.Range("K31").Value = 1864.5
$ws.Range("L31").Value = 9493.394
$ws.Range("M31").Value = -1569.5
$ws.Range("N31").Value = -10083.394
$ws.Range("H34").Value = 7002.3267
$ws.Range("I34").Value = 1864.5
$ws.Range("J34").Value = 9493.394
$ws.Range("K34").Value = 1864.5
$ws.Range("L34").Value = 9493.394
$ws.Range("M34").Value = -1662.5
$ws.Range("N34").Value = -9897.394
$ws.Range("H134").Value = 3965
$ws.Range("I134").Value = 2960.0715
$ws.Range("K134").Value = 8880.2145
$ws.Range("M134").Value = -6345.2145

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H109").Value = 3864.6956
$ws.Range("I109").Value = 2431.6667
$ws.Range("J109").Value = 4785.9287
$ws.Range("K109").Value = 7295.000100000001
$ws.Range("L109").Value = 14357.7861
$ws.Range("M109").Value = -6255.000100000001
$ws.Range("N109").Value = -16437.7861
$ws.Range("H112").Value = 5174.6
$ws.Range("I112").Value = 3991
$ws.Range("J112").Value = 6950
$ws.Range("K112").Value = 11973
$ws.Range("L112").Value = 20850
$ws.Range("M112").Value = -10865
$ws.Range("N112").Value = -23066
$ws.Range("H123").Value = 4870.4
$ws.Range("I123").Value = 3838
$ws.Range("J123").Value = 9000
$ws.Range("K123").Value = 11514
$ws.Range("L123").Value = 27000
$ws.Range("M123").Value = -9064
$ws.Range("N123").Value = -31900
$ws.Range("H131").Value = 15200.5
$ws.Range("I131").Value = 10346
$ws.Range("K131").Value = 31038
$ws.Range("M131").Value = -25998
$ws.Range("H139").Value = 951.4666999999999
$ws.Range("I139").Value = 790.1539
$ws.Range("K139").Value = 2370.4617
$ws.Range("M139").Value = 2769.5383

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 984.1818
$ws.Range("I97").Value = 380.94736
$ws.Range("K97").Value = 380.94736
$ws.Range("M97").Value = 115.05264
$ws.Range("H134").Value = 50332.168
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 50332.168
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 150996.504
$ws.Range("N134").Value = -156066.504
$ws.Range("H136").Value = 29315.074
$ws.Range("J136").Value = 29315.074
$ws.Range("L136").Value = 87945.22200000001
$ws.Range("N136").Value = -93045.22200000001
$ws.Range("M134").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 55267.844
$ws.Range("I22").Value = 143928.28
$ws.Range("J22").Value = 3549.25
$ws.Range("K22").Value = 143928.28
$ws.Range("L22").Value = 3549.25
$ws.Range("M22").Value = -143633.28
$ws.Range("N22").Value = -4139.25
$ws.Range("H25").Value = 47500
$ws.Range("I25").Value = 47500
$ws.Range("K25").Value = 47500
$ws.Range("M25").Value = -47270
$ws.Range("H27").Value = 55267.844
$ws.Range("I27").Value = 143928.28
$ws.Range("J27").Value = 3549.25
$ws.Range("K27").Value = 143928.28
$ws.Range("L27").Value = 3549.25
$ws.Range("M27").Value = -143821.28
$ws.Range("N27").Value = -3763.25
$ws.Range("H42").Value = 18316.666
$ws.Range("I42").Value = 15000
$ws.Range("K42").Value = 15000
$ws.Range("M42").Value = -14437
$ws.Range("H43").Value = 259500
$ws.Range("I43").Value = 5000
$ws.Range("J43").Value = 386750
$ws.Range("K43").Value = 5000
$ws.Range("L43").Value = 386750
$ws.Range("M43").Value = -4807
$ws.Range("N43").Value = -387136
$ws.Range("H49").Value = 18316.666
$ws.Range("I49").Value = 15000
$ws.Range("K49").Value = 15000
$ws.Range("M49").Value = -14853
$ws.Range("H68").Value = 6660.6665
$ws.Range("J68").Value = 6368.3125
$ws.Range("L68").Value = 6368.3125
$ws.Range("N68").Value = -7866.3125
$ws.Range("H71").Value = 6660.6665
$ws.Range("J71").Value = 6368.3125
$ws.Range("L71").Value = 31841.5625
$ws.Range("N71").Value = -39329.5625
$ws.Range("H138").Value = 103013.25
$ws.Range("J138").Value = 103013.25
$ws.Range("L138").Value = 103013.25
$ws.Range("N138").Value = -113293.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 73050.75
$ws.Range("I46").Value = 45000
$ws.Range("J46").Value = 77058
$ws.Range("K46").Value = 45000
$ws.Range("L46").Value = 77058
$ws.Range("M46").Value = -44769
$ws.Range("N46").Value = -77520
$ws.Range("H81").Value = 842.8570999999999
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("H84").Value = 842.8570999999999
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("H132").Value = 3413599.8
$ws.Range("I132").Value = 4474180
$ws.Range("J132").Value = 4592.5713
$ws.Range("K132").Value = 13422540
$ws.Range("L132").Value = 13777.7139
$ws.Range("M132").Value = -13420010
$ws.Range("N132").Value = -18837.7139
$ws.Range("H134").Value = 73050.75
$ws.Range("I134").Value = 45000
$ws.Range("J134").Value = 77058
$ws.Range("K134").Value = 135000
$ws.Range("L134").Value = 231174
$ws.Range("M134").Value = -132465
$ws.Range("N134").Value = -236244
$ws.Range("N81").ClearContents()
$ws.Range("N84").ClearContents()
